$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "52 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "2|    |"
$t.Cell(1,2).Range.Text = "22 x 23" + [char]11 + "  2    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"
$t.Cell(1,3).Range.Text = "81 x 75" + [char]11 + "  7    5" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"
$t.Cell(2,1).Range.Text = "22 x 63" + [char]11 + "  6    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"
$t.Cell(2,2).Range.Text = "91 x 98" + [char]11 + "  9    8" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "1|    |"
$t.Cell(2,3).Range.Text = "60 x 83" + [char]11 + "  8    3" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "0|    |"
$t.Cell(3,1).Range.Text = "50 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "0|    |"
$t.Cell(3,2).Range.Text = "15 x 38" + [char]11 + "  3    8" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "5|    |"
$t.Cell(3,3).Range.Text = "51 x 18" + [char]11 + "  1    8" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "1|    |"
$t.Cell(4,1).Range.Text = "52 x 37" + [char]11 + "  3    7" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "2|    |"
$t.Cell(4,2).Range.Text = "50 x 80" + [char]11 + "  8    0" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "0|    |"
$t.Cell(4,3).Range.Text = "77 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "7|    |"
$t.Cell(5,1).Range.Text = "14 x 97" + [char]11 + "  9    7" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "4|    |"
$t.Cell(5,2).Range.Text = "74 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "4|    |"
$t.Cell(5,3).Range.Text = "63 x 58" + [char]11 + "  5    8" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "3|    |"
